# Author commit: swap the deck's applied colour theme ("Integral" -> the
# previous "Office Theme" palette) and pick a different built-in table
# style for the cash-flow table on slide 16.

$p = $ppt.ActivePresentation

# --- 1. Re-colour the presentation's live theme (clrScheme) -------------
# The 12 theme colour slots, in MsoThemeColorSchemeIndex order
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink), switched from the
# "Integral" palette back to the stock "Office Theme" palette.
$officeThemeRGB = @(
    0,            # Dark 1    - 000000
    16777215,     # Light 1   - FFFFFF
    6968388,      # Dark 2    - 44546A
    15132391,     # Light 2   - E7E6E6
    13998939,     # Accent 1  - 5B9BD5
    3243501,      # Accent 2  - ED7D31
    10855845,     # Accent 3  - A5A5A5
    49407,        # Accent 4  - FFC000
    12874308,     # Accent 5  - 4472C4
    4697456,      # Accent 6  - 70AD47
    12673797,     # Hyperlink - 0563C1
    7491477       # Followed Hyperlink - 954F72
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeRGB[$i - 1]
}

# --- 2. Switch the slide-16 table to the new table style ----------------
$tableSlide = $p.Slides.Item(16)
$tableShape = $tableSlide.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{D9559F75-DD18-46FA-B1D3-DCC17715DF81}", $false)
